$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto data (prices, 1h volume %, and two reordered-row swaps)

$ws.Range("D2").Value = "'36.926.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.38%  "

$ws.Range("D3").Value = "'1.985.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.75%  "

$ws.Range("D4").Value = "'0.996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.43%  "

$ws.Range("D5").Value = "'264.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.59%  "

$ws.Range("D6").Value = "'0.609"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.45%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("D8").Value = "'55.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -7.45%  "

$ws.Range("D9").Value = "'0.373"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.25%  "

$ws.Range("D10").Value = "'0.0758"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.08%  "

$ws.Range("D11").Value = "'0.101"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.50%  "

$ws.Range("D12").Value = "'14.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.46%  "

$ws.Range("D13").Value = "'2.255.09"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.50%  "

$ws.Range("D14").Value = "'21.30"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.61%  "

$ws.Range("D15").Value = "'0.769"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -9.62%  "

$ws.Range("D16").Value = "'5.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.47%  "

$ws.Range("D17").Value = "'1.987.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.50%  "

$ws.Range("D18").Value = "'36.630.96"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.90%  "

$ws.Range("D19").Value = "'69.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.77%  "

$ws.Range("D20").Value = "'0.0₃0825"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.41%  "

$ws.Range("D21").Value = "'233.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.54%  "

$ws.Range("D22").Value = "'5.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.64%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "'2.56"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.12%  "

$ws.Range("D25").Value = "'2.37"
$ws.Range("D25").Style = "Normal"

$ws.Range("D26").Value = "'163.52"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.62%  "

$ws.Range("D27").Value = "'8.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.78%  "

$ws.Range("D28").Value = "'19.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.43%  "

$ws.Range("B29").Value = "ImmutableX"
$ws.Range("C29").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D29").Value = "'1.33"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.12%  "

$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").Value = "'0.122"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -12.13%  "

$ws.Range("D31").Value = "'0.118"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.88%  "

$ws.Range("D32").Value = "'4.52"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.01%  "

$ws.Range("D33").Value = "'0.0624"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -8.17%  "

$ws.Range("D34").Value = "'4.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.70%  "

$ws.Range("D35").Value = "'2.38"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -6.27%  "

$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "'3.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.64%  "

$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "'1.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.45%  "

$ws.Range("B38").Value = "BinanceUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D38").Value = "'0.998"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.27%  "

$ws.Range("E39").Value = "  -3.44%  "

$ws.Range("D40").Value = "'2.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.24%  "

$ws.Range("D41").Value = "'1.444.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.87%  "

$ws.Range("D42").Value = "'1.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.43%  "

$ws.Range("D43").Value = "'0.0910"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.33%  "

$ws.Range("D44").Value = "'0.0205"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.27%  "

$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'88.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.02%  "

$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'15.41"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -7.77%  "

$ws.Range("D47").Value = "'1.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.88%  "

$ws.Range("E48").Value = "  -0.19%  "

$ws.Range("D49").Value = "'6.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.48%  "

$ws.Range("D50").Value = "'2.146.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.65%  "

$ws.Range("D51").Value = "'1.92"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -9.50%  "
